$wb = $excel.ActiveWorkbook

# This script applies a market-data refresh to the profit-calculation
# columns (H: currentAveragePrice, I/J: NQ/HQ average price, K/L: Leve
# price NQ/HQ, M/N: Leve profit NQ/HQ) across several sheets, matching a
# scheduled-runner data sync commit.

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 2307.1667
$ws.Range("J4").Value = 4998.5
$ws.Range("L4").Value = 4998.5
$ws.Range("N4").Value = -5226.5

# Row 64: Forged from the Void
$ws.Range("H64").Value = 3125
$ws.Range("I64").Value = 3125
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3125
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2877
$ws.Range("N64").ClearContents()

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3125
$ws.Range("I67").Value = 3125
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3125
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2267
$ws.Range("N67").ClearContents()

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 107079.6
$ws.Range("I86").Value = 8966.166999999999
$ws.Range("J86").Value = 254249.75
$ws.Range("K86").Value = 8966.166999999999
$ws.Range("L86").Value = 254249.75
$ws.Range("M86").Value = -7843.166999999999
$ws.Range("N86").Value = -256495.75

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 107079.6
$ws.Range("I89").Value = 8966.166999999999
$ws.Range("J89").Value = 254249.75
$ws.Range("K89").Value = 44830.835
$ws.Range("L89").Value = 1271248.75
$ws.Range("M89").Value = -39214.835
$ws.Range("N89").Value = -1282480.75

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 2096.3333
$ws.Range("I111").Value = 290
$ws.Range("K111").Value = 870
$ws.Range("M111").Value = 2197

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3172.3333
$ws.Range("I138").Value = 2508.92
$ws.Range("J138").Value = 4680.091
$ws.Range("K138").Value = 7526.76
$ws.Range("L138").Value = 14040.273
$ws.Range("M138").Value = -2386.76
$ws.Range("N138").Value = -24320.273

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3144.027
$ws.Range("I61").Value = 3150.9688
$ws.Range("K61").Value = 3150.9688
$ws.Range("M61").Value = -2938.9688

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 8545.362999999999
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 8545.362999999999
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1532.2188
$ws.Range("I74").Value = 1197.4445
$ws.Range("J74").Value = 3340
$ws.Range("K74").Value = 1197.4445
$ws.Range("L74").Value = 3340
$ws.Range("M74").Value = -323.4445000000001
$ws.Range("N74").Value = -5088

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1532.2188
$ws.Range("I77").Value = 1197.4445
$ws.Range("J77").Value = 3340
$ws.Range("K77").Value = 5987.2225
$ws.Range("L77").Value = 16700
$ws.Range("M77").Value = -1619.2225
$ws.Range("N77").Value = -25436

# Row 112: Wrapped Knuckles
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3144.027
$ws.Range("I136").Value = 3150.9688
$ws.Range("K136").Value = 9452.9064
$ws.Range("M136").Value = -6902.9064

$ws = $wb.Worksheets.Item("BSM")
# Row 25: Tools of the Trade
$ws.Range("H25").Value = 6682.6
$ws.Range("I25").Value = 5353.25
$ws.Range("K25").Value = 5353.25
$ws.Range("M25").Value = -5118.25

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4979.2856
$ws.Range("I134").Value = 4760
$ws.Range("J134").Value = 5527.5
$ws.Range("K134").Value = 14280
$ws.Range("L134").Value = 16582.5
$ws.Range("M134").Value = -11745
$ws.Range("N134").Value = -21652.5

$ws = $wb.Worksheets.Item("CRP")
# Row 36: Toys of Summer
$ws.Range("H36").Value = 4999.5
$ws.Range("I36").Value = 4999.5
$ws.Range("K36").Value = 4999.5
$ws.Range("M36").Value = -4611.5

# Row 40: Ceremonial Spears
$ws.Range("H40").Value = 4999.5
$ws.Range("I40").Value = 4999.5
$ws.Range("K40").Value = 4999.5
$ws.Range("M40").Value = -4839.5

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2982.0967
$ws.Range("I58").Value = 2265.5
$ws.Range("K58").Value = 2265.5
$ws.Range("M58").Value = -2062.5

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3557.5386
$ws.Range("I134").Value = 3624.9
$ws.Range("K134").Value = 10874.7
$ws.Range("M134").Value = -8339.700000000001

# Row 136: Turali Quality
$ws.Range("H136").Value = 2982.0967
$ws.Range("I136").Value = 2265.5
$ws.Range("K136").Value = 6796.5
$ws.Range("M136").Value = -4246.5

$ws = $wb.Worksheets.Item("CUL")
# Row 60: Drinking to Your Health
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1438.9697
$ws.Range("J131").Value = 1445.6495
$ws.Range("L131").Value = 4336.9485
$ws.Range("N131").Value = -14416.9485

# Row 132: More Mezcal
$ws.Range("H132").Value = 3910.5625
$ws.Range("I132").Value = 2698.6667
$ws.Range("K132").Value = 24288.0003
$ws.Range("M132").Value = -21758.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 1998.3334
$ws.Range("I80").Value = 1998.3334
$ws.Range("K80").Value = 1998.3334
$ws.Range("M80").Value = -1000.3334

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1998.3334
$ws.Range("I83").Value = 1998.3334
$ws.Range("K83").Value = 9991.666999999999
$ws.Range("M83").Value = -4999.666999999999

# Row 99: Needle in a Hingan Stack
$ws.Range("H99").Value = 300495.72
$ws.Range("I99").Value = 300495.72
$ws.Range("K99").Value = 300495.72
$ws.Range("M99").Value = -298249.72

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 3318.889
$ws.Range("I7").Value = 3278.8235
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 3278.8235
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -3166.8235
$ws.Range("N7").Value = -4224

# Row 16: Saddle Sore
$ws.Range("H16").Value = 13198.4
$ws.Range("I16").Value = 15998
$ws.Range("K16").Value = 15998
$ws.Range("M16").Value = -15828

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1013.7
$ws.Range("J22").Value = 1267.6
$ws.Range("L22").Value = 1267.6
$ws.Range("N22").Value = -1857.6

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1013.7
$ws.Range("J27").Value = 1267.6
$ws.Range("L27").Value = 1267.6
$ws.Range("N27").Value = -1481.6

# Row 42: Slave to Fashion
$ws.Range("H42").Value = 844666.5
$ws.Range("J42").Value = 1011399.8
$ws.Range("L42").Value = 1011399.8
$ws.Range("N42").Value = -1012525.8

# Row 43: Subordinate Clause
$ws.Range("H43").Value = 1260587.4
$ws.Range("I43").Value = 10012
$ws.Range("J43").Value = 1885875
$ws.Range("K43").Value = 10012
$ws.Range("L43").Value = 1885875
$ws.Range("M43").Value = -9819
$ws.Range("N43").Value = -1886261

# Row 45: Soft Shoe Shuffle
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29593

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3516.6667

# Row 49: First They Came for the Heretics
$ws.Range("H49").Value = 844666.5
$ws.Range("J49").Value = 1011399.8
$ws.Range("L49").Value = 1011399.8
$ws.Range("N49").Value = -1011693.8

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 10103359
$ws.Range("I61").Value = 13891019
$ws.Range("J61").Value = 2933
$ws.Range("K61").Value = 13891019
$ws.Range("L61").Value = 2933
$ws.Range("M61").Value = -13890817
$ws.Range("N61").Value = -3337

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 3598.8
$ws.Range("I68").Value = 3598.8
$ws.Range("K68").Value = 3598.8
$ws.Range("M68").Value = -2849.8

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 3598.8
$ws.Range("I71").Value = 3598.8
$ws.Range("K71").Value = 17994
$ws.Range("M71").Value = -14250

# Row 113: Peace in Rest
$ws.Range("H113").Value = 10103359
$ws.Range("I113").Value = 13891019
$ws.Range("J113").Value = 2933
$ws.Range("K113").Value = 13891019
$ws.Range("L113").Value = 2933
$ws.Range("M113").Value = -13888849
$ws.Range("N113").Value = -7273

# Row 119: Fit for a Friend
$ws.Range("H119").Value = 55000
$ws.Range("J119").Value = 55000
$ws.Range("L119").Value = 55000
$ws.Range("N119").Value = -64676

# Row 126: Battered Books
$ws.Range("H126").Value = 3318.889
$ws.Range("I126").Value = 3278.8235
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 9836.470499999999
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -7366.470499999999
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("WVR")
# Row 25: A Drag of a Doublet
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 30: The Telltale Tress
$ws.Range("H30").Value = 11990
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 11990
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 11990
$ws.Range("N30").Value = -12204
$ws.Range("M30").ClearContents()

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7145.4707
$ws.Range("I132").Value = 4284.5
$ws.Range("K132").Value = 12853.5
$ws.Range("M132").Value = -10323.5
